$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172093272209167
$ws.Range("B1").Value = 4.173731327056885
$ws.Range("C1").Value = 3.794234752655029
$ws.Range("D1").Value = 1.803316950798035
$ws.Range("E1").Value = 1.302623629570007
